$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11.20414814319263
$ws.Range("D2").Value = 8.023658105252926
$ws.Range("E2").Value = 12.79246775424638
$ws.Range("F2").Value = 31.14583031427596
$ws.Range("G2").Value = 33.31539252836562
$ws.Range("H2").Value = 15.52691948833641
$ws.Range("J2").Value = 9.761376319509955
$ws.Range("L2").Value = 9.334405276761945
$ws.Range("M2").Value = 60.06474090040093
$ws.Range("O2").Value = 24.22361037799299
$ws.Range("C3").Value = 11.39707426622343
$ws.Range("D3").Value = 8.097584704228877
$ws.Range("E3").Value = 12.85661749350699
$ws.Range("F3").Value = 31.57657441623798
$ws.Range("G3").Value = 33.78608575261823
$ws.Range("H3").Value = 15.67401016675538
$ws.Range("J3").Value = 9.809711319168132
$ws.Range("L3").Value = 9.319941116971309
$ws.Range("M3").Value = 56.68553178440989
$ws.Range("O3").Value = 24.50321391722862
$ws.Range("C4").Value = 11.51994177451187
$ws.Range("D4").Value = 8.1448388621285
$ws.Range("E4").Value = 12.90128259781346
$ws.Range("F4").Value = 31.85473673523169
$ws.Range("G4").Value = 34.09571533214463
$ws.Range("H4").Value = 15.76893405266173
$ws.Range("J4").Value = 9.842480319873333
$ws.Range("L4").Value = 9.313369814829951
$ws.Range("M4").Value = 54.49502894174223
$ws.Range("O4").Value = 24.68466313253831
$ws.Range("C5").Value = 11.57112780740304
$ws.Range("D5").Value = 8.164565721696498
$ws.Range("E5").Value = 12.92079071278689
$ws.Range("F5").Value = 31.9715078167225
$ws.Range("G5").Value = 34.22698436554189
$ws.Range("H5").Value = 15.80877214921539
$ws.Range("J5").Value = 9.856602729104846
$ws.Range("L5").Value = 9.311271630682075
$ws.Range("M5").Value = 53.57359655211093
$ws.Range("O5").Value = 24.7610427891878
$ws.Range("C6").Value = 11.57969486192313
$ws.Range("D6").Value = 8.167869820274408
$ws.Range("E6").Value = 12.92410836281738
$ws.Range("F6").Value = 31.99110333687814
$ws.Range("G6").Value = 34.24908642474099
$ws.Range("H6").Value = 15.81545695253328
$ws.Range("J6").Value = 9.858993933993077
$ws.Range("L6").Value = 9.310958202948225
$ws.Range("M6").Value = 53.41886529053974
$ws.Range("O6").Value = 24.7738722487251
$ws.Range("C7").Value = 11.52062755658078
$ws.Range("D7").Value = 8.145102998192659
$ws.Range("E7").Value = 12.90154042647435
$ws.Range("F7").Value = 31.85629774220846
$ws.Range("G7").Value = 34.09746517263132
$ws.Range("H7").Value = 15.7694666473553
$ws.Range("J7").Value = 9.842667677864664
$ws.Range("L7").Value = 9.313339172522634
$ws.Range("M7").Value = 54.48271827461922
$ws.Range("O7").Value = 24.68568337038488
$ws.Range("C8").Value = 11.2697584096562
$ws.Range("D8").Value = 8.048762536205411
$ws.Range("E8").Value = 12.81347905067132
$ws.Range("F8").Value = 31.29149642517997
$ws.Range("G8").Value = 33.47335005892062
$ws.Range("H8").Value = 15.57667761163953
$ws.Range("J8").Value = 9.777395446170029
$ws.Range("L8").Value = 9.328938012267173
$ws.Range("M8").Value = 58.92374672479591
$ws.Range("O8").Value = 24.31797758820882
$ws.Range("C9").Value = 10.81243490247374
$ws.Range("D9").Value = 7.874528925666216
$ws.Range("E9").Value = 12.68359261262298
$ws.Range("F9").Value = 30.29360587845179
$ws.Range("G9").Value = 32.41748880132307
$ws.Range("H9").Value = 15.23534527451741
$ws.Range("J9").Value = 9.674321894833009
$ws.Range("L9").Value = 9.377901194080787
$ws.Range("M9").Value = 66.70580815806292
$ws.Range("O9").Value = 23.67535742620534
$ws.Range("C10").Value = 10.49701979520897
$ws.Range("D10").Value = 7.755343327099602
$ws.Range("E10").Value = 12.61557513792997
$ws.Range("F10").Value = 29.62897060590287
$ws.Range("G10").Value = 31.75059326612092
$ws.Range("H10").Value = 15.00717455738372
$ws.Range("J10").Value = 9.614361794295798
$ws.Range("L10").Value = 9.425135823563847
$ws.Range("M10").Value = 71.85057707223395
$ws.Range("O10").Value = 23.25238542123691
$ws.Range("C11").Value = 10.35787775428549
$ws.Range("D11").Value = 7.703010427678046
$ws.Range("E11").Value = 12.59087794084521
$ws.Range("F11").Value = 29.34189519385657
$ws.Range("G11").Value = 31.47234623074625
$ws.Range("H11").Value = 14.90833543794445
$ws.Range("J11").Value = 9.590639462346747
$ws.Range("L11").Value = 9.449079015511797
$ws.Range("M11").Value = 74.06606738367849
$ws.Range("O11").Value = 23.07095783487689
$ws.Range("C12").Value = 10.30580282830303
$ws.Range("D12").Value = 7.683462291298127
$ws.Range("E12").Value = 12.582448544228
$ws.Range("F12").Value = 29.2354218610498
$ws.Range("G12").Value = 31.37072898912849
$ws.Range("H12").Value = 14.87162565662956
$ws.Range("J12").Value = 9.582178853162594
$ws.Range("L12").Value = 9.458499181790277
$ws.Range("M12").Value = 74.88704610621501
$ws.Range("O12").Value = 23.00386479723084
$ws.Range("C13").Value = 10.31699088671566
$ws.Range("D13").Value = 7.687660389609589
$ws.Range("E13").Value = 12.58422251677146
$ws.Range("F13").Value = 29.25825268152018
$ws.Range("G13").Value = 31.3924451101302
$ws.Range("H13").Value = 14.87949971692834
$ws.Range("J13").Value = 9.583977569956712
$ws.Range("L13").Value = 9.456454660714048
$ws.Range("M13").Value = 74.7110333033356
$ws.Range("O13").Value = 23.01824239293949
$ws.Range("C14").Value = 10.35358124432147
$ws.Range("D14").Value = 7.701396806448684
$ws.Range("E14").Value = 12.59016580841146
$ws.Range("F14").Value = 29.33309051642152
$ws.Range("G14").Value = 31.4639101832467
$ws.Range("H14").Value = 14.9053008694523
$ws.Range("J14").Value = 9.589932863543007
$ws.Range("L14").Value = 9.449846934828251
$ws.Range("M14").Value = 74.1339701438085
$ws.Range("O14").Value = 23.06540561100907
$ws.Range("C15").Value = 10.37607372601351
$ws.Range("D15").Value = 7.709845759785924
$ws.Range("E15").Value = 12.59392718381766
$ws.Range("F15").Value = 29.37922327798881
$ws.Range("G15").Value = 31.50817696832508
$ws.Range("H15").Value = 14.92119855879949
$ws.Range("J15").Value = 9.593649048833177
$ws.Range("L15").Value = 9.445845539368184
$ws.Range("M15").Value = 73.77816074696398
$ws.Range("O15").Value = 23.09450499538294
$ws.Range("C16").Value = 10.50619964428623
$ws.Range("D16").Value = 7.75880117818213
$ws.Range("E16").Value = 12.61731693265901
$ws.Range("F16").Value = 29.64804228265539
$ws.Range("G16").Value = 31.76929441247578
$ws.Range("H16").Value = 15.01373414488174
$ws.Range("J16").Value = 9.615984588841203
$ws.Range("L16").Value = 9.423620523227692
$ws.Range("M16").Value = 71.70327403891737
$ws.Range("O16").Value = 23.26446570414012
$ws.Range("C17").Value = 10.58713328474275
$ws.Range("D17").Value = 7.789315238278966
$ws.Range("E17").Value = 12.63328271187635
$ws.Range("F17").Value = 29.81689256728783
$ws.Range("G17").Value = 31.93601134953485
$ws.Range("H17").Value = 15.07177519331852
$ws.Range("J17").Value = 9.630604978156612
$ws.Range("L17").Value = 9.410615437653231
$ws.Range("M17").Value = 70.39838668958697
$ws.Range("O17").Value = 23.37156527686202
$ws.Range("C18").Value = 10.63409342232396
$ws.Range("D18").Value = 7.80704367764432
$ws.Range("E18").Value = 12.64305202304521
$ws.Range("F18").Value = 29.91544552125617
$ws.Range("G18").Value = 32.03425904352786
$ws.Range("H18").Value = 15.1056250710519
$ws.Range("J18").Value = 9.639348028575631
$ws.Range("L18").Value = 9.403366415809954
$ws.Range("M18").Value = 69.63609378426828
$ws.Range("O18").Value = 23.43419811376165
$ws.Range("C19").Value = 10.6500638819606
$ws.Range("D19").Value = 7.813076776481825
$ws.Range("E19").Value = 12.64645967402903
$ws.Range("F19").Value = 29.94905911158967
$ws.Range("G19").Value = 32.06792474538113
$ws.Range("H19").Value = 15.11716596319491
$ws.Range("J19").Value = 9.642365260412637
$ws.Range("L19").Value = 9.400951713294917
$ws.Range("M19").Value = 69.37597647230463
$ws.Range("O19").Value = 23.45558087723052
$ws.Range("C20").Value = 10.57847547065793
$ws.Range("D20").Value = 7.786048601615897
$ws.Range("E20").Value = 12.63152227685364
$ws.Range("F20").Value = 29.79876939772978
$ws.Range("G20").Value = 31.91801919691418
$ws.Range("H20").Value = 15.06554833496932
$ws.Range("J20").Value = 9.629013984987795
$ws.Range("L20").Value = 9.411975924874485
$ws.Range("M20").Value = 70.53851113867728
$ws.Range("O20").Value = 23.36005732530258
$ws.Range("C21").Value = 10.34281714942248
$ws.Range("D21").Value = 7.697354796912818
$ws.Range("E21").Value = 12.58839487097022
$ws.Range("F21").Value = 29.31104781221933
$ws.Range("G21").Value = 31.44281629813695
$ws.Range("H21").Value = 14.89770289976253
$ws.Range("J21").Value = 9.588169374646281
$ws.Range("L21").Value = 9.451778191757583
$ws.Range("M21").Value = 74.30395555523161
$ws.Range("O21").Value = 23.05150869523355
$ws.Range("C22").Value = 10.19238184567212
$ws.Range("D22").Value = 7.640956489691884
$ws.Range("E22").Value = 12.56559988442284
$ws.Range("F22").Value = 29.00534720985054
$ws.Range("G22").Value = 31.15416322600845
$ws.Range("H22").Value = 14.79219681768065
$ws.Range("J22").Value = 9.564526355538892
$ws.Range("L22").Value = 9.47985116159709
$ws.Range("M22").Value = 76.66010355878564
$ws.Range("O22").Value = 22.85925203918506
$ws.Range("C23").Value = 10.27234746081052
$ws.Range("D23").Value = 7.670914456221991
$ws.Range("E23").Value = 12.57726432164984
$ws.Range("F23").Value = 29.16729681511723
$ws.Range("G23").Value = 31.30617043971471
$ws.Range("H23").Value = 14.8481219728791
$ws.Range("J23").Value = 9.576861980906324
$ws.Range("L23").Value = 9.464679568751453
$ws.Range("M23").Value = 75.41217054025105
$ws.Range("O23").Value = 22.96099267060961
$ws.Range("C24").Value = 10.58238832739492
$ws.Range("D24").Value = 7.787524870117589
$ws.Range("E24").Value = 12.63231633253661
$ws.Range("F24").Value = 29.80695827992677
$ws.Range("G24").Value = 31.92614598412193
$ws.Range("H24").Value = 15.06836199889179
$ws.Range("J24").Value = 9.629732222335052
$ws.Range("L24").Value = 9.411360139153201
$ws.Range("M24").Value = 70.47519854907537
$ws.Range("O24").Value = 23.36525677121083
$ws.Range("C25").Value = 10.93249946368272
$ws.Range("D25").Value = 7.920104685381524
$ws.Range("E25").Value = 12.71401423256485
$ws.Range("F25").Value = 30.55165828297515
$ws.Range("G25").Value = 32.68449983413585
$ws.Range("H25").Value = 15.32373168312031
$ws.Range("J25").Value = 9.699479920362496
$ws.Range("L25").Value = 9.362681346232909
$ws.Range("M25").Value = 64.70077377476073
$ws.Range("O25").Value = 23.84067225498913
